$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Aple"
$ws.Range("C2").Value = "Real"
$ws.Range("D2").Value = "COButts"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Mommy"
$ws.Range("C3").Value = "Fake"
$ws.Range("D3").Value = "Butts"
$ws.Range("E3").Value = "2021-08-12T16:04:45.663Z"
